$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "{'batch_size': 16, 'epochs': 100, 'layers_struct': [{'units': 100, 'dropout': 0.3}, {'units': 50, 'dropout': 0.2}]}"

$ws.Range("I2").Value = "rmse"
$ws.Range("J2").Value = 0.08602705299854296

$ws.Range("N2").Value = 1.788979776925947
$ws.Range("P2").Value = 1.403455132887428
$ws.Range("Q2").Value = 0.9436809236165411
$ws.Range("R2").Value = 2.7042386776883
$ws.Range("T2").Value = 2.196291687818989
$ws.Range("U2").Value = 0.8899396914830616
